# Applies the CryCompanywiseStockReport_1 batch-reconciliation edit:
# for several products that were stocked/issued in multiple batches, the
# batch no. (B), rate (D/E), qty (F) and value (G = D*F) on the affected
# rows are corrected/re-paired, and every "Sub Total:" / "Grand Total:" row
# that rolls those rows up is recomputed to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B313").Value = 57854
$ws.Range("F313").Value = 2
$ws.Range("G313").Value = 611.6799999999999
$ws.Range("B314").Value = 62997
$ws.Range("F314").Value = 0
$ws.Range("G314").Value = 0

$ws.Range("B346").Value = 63520
$ws.Range("E346").Value = 153.4
$ws.Range("F346").Value = 87
$ws.Range("G346").Value = 12552.36
$ws.Range("B347").Value = 55373
$ws.Range("E347").Value = 163.62
$ws.Range("F347").Value = -94
$ws.Range("G347").Value = -13562.32

$ws.Range("B350").Value = 63571
$ws.Range("E350").Value = 152.53
$ws.Range("F350").Value = 15
$ws.Range("G350").Value = 2152.2
$ws.Range("B351").Value = 63531
$ws.Range("F351").Value = 80
$ws.Range("G351").Value = 11478.4
$ws.Range("B352").Value = 57802
$ws.Range("E352").Value = 162.71
$ws.Range("F352").Value = -79
$ws.Range("G352").Value = -11334.92

$ws.Range("B372").Value = 57885
$ws.Range("E372").Value = 62.28
$ws.Range("F372").Value = 4
$ws.Range("G372").Value = 208.52
$ws.Range("B373").Value = 63652
$ws.Range("E373").Value = 55.42
$ws.Range("F373").Value = 165
$ws.Range("G373").Value = 8601.450000000001

$ws.Range("B375").Value = 61605
$ws.Range("E375").Value = 133.78
$ws.Range("F375").Value = -13
$ws.Range("G375").Value = -1455.48
$ws.Range("B376").Value = 63563
$ws.Range("E376").Value = 119.04
$ws.Range("F376").Value = 2
$ws.Range("G376").Value = 223.92

$ws.Range("B389").Value = 62865
$ws.Range("F389").Value = 13
$ws.Range("G389").Value = 1037.53
$ws.Range("B390").Value = 57817
$ws.Range("F390").Value = 3
$ws.Range("G390").Value = 239.43

$ws.Range("B400").Value = 57835
$ws.Range("F400").Value = 1
$ws.Range("G400").Value = 59.13
$ws.Range("B401").Value = 62933
$ws.Range("F401").Value = 116
$ws.Range("G401").Value = 6859.08

$ws.Range("B419").Value = 63007
$ws.Range("F419").Value = 838
$ws.Range("G419").Value = 143574.54
$ws.Range("B420").Value = 57856
$ws.Range("F420").Value = 2
$ws.Range("G420").Value = 342.66
$ws.Range("B421").Value = 63008
$ws.Range("F421").Value = 433
$ws.Range("G421").Value = 65456.61
$ws.Range("B422").Value = 57857
$ws.Range("F422").Value = 3
$ws.Range("G422").Value = 453.51

$ws.Range("B434").Value = 483083.13

$ws.Range("F536").Value = 44
$ws.Range("G536").Value = 4643.76

$ws.Range("B546").Value = 72196.71000000001

$ws.Range("B583").Value = 65066
$ws.Range("E583").Value = 13.61
$ws.Range("F583").Value = 231
$ws.Range("G583").Value = 2959.11
$ws.Range("B584").Value = 53263
$ws.Range("E584").Value = 15.29
$ws.Range("F584").Value = -309
$ws.Range("G584").Value = -3958.29

$ws.Range("F589").Value = 56
$ws.Range("G589").Value = 717.36

$ws.Range("B599").Value = 64925
$ws.Range("E599").Value = 13.97
$ws.Range("F599").Value = 269
$ws.Range("G599").Value = 3537.35
$ws.Range("B600").Value = 45709
$ws.Range("E600").Value = 15.69
$ws.Range("F600").Value = -300
$ws.Range("G600").Value = -3945

$ws.Range("B606").Value = 11703.9

$ws.Range("F609").Value = 86
$ws.Range("G609").Value = 2691.8

$ws.Range("B625").Value = 33479.81

$ws.Range("B709").Value = 64833
$ws.Range("E709").Value = 34.9
$ws.Range("F709").Value = 97
$ws.Range("G709").Value = 3184.51
$ws.Range("B710").Value = 60025
$ws.Range("E710").Value = 37.22
$ws.Range("F710").Value = -98
$ws.Range("G710").Value = -3217.34

$ws.Range("B715").Value = 60031
$ws.Range("E715").Value = 111.69
$ws.Range("F715").Value = -5
$ws.Range("G715").Value = -492.5
$ws.Range("B716").Value = 64836
$ws.Range("E716").Value = 104.71
$ws.Range("F716").Value = 6
$ws.Range("G716").Value = 591

$ws.Range("B720").Value = 64830
$ws.Range("E720").Value = 34.9
$ws.Range("F720").Value = 114
$ws.Range("G720").Value = 3742.62
$ws.Range("B721").Value = 60022
$ws.Range("E721").Value = 37.22
$ws.Range("F721").Value = -113
$ws.Range("G721").Value = -3709.79

$ws.Range("F855").Value = 123
$ws.Range("G855").Value = 10031.88

$ws.Range("B859").Value = 63150
$ws.Range("D859").Value = 75.68000000000001
$ws.Range("E859").Value = 80.45
$ws.Range("F859").Value = 151
$ws.Range("G859").Value = 11427.68
$ws.Range("B860").Value = 61428
$ws.Range("D860").Value = 69.16
$ws.Range("E860").Value = 73.52
$ws.Range("F860").Value = 1
$ws.Range("G860").Value = 69.16

$ws.Range("B870").Value = 302528.9

$ws.Range("B962").Value = 4363862.4
$ws.Range("B963").Value = 4363862.4
